# Generate Report for Handoff
# Replaces the e2e test-run artifact data (file names, statuses, timestamps)
# with the values from the latest handoff run across all three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: remove a single hyperlink (by its top-left cell address) from a
# worksheet without disturbing any other hyperlinks on that sheet.
# ---------------------------------------------------------------------------
function Remove-HyperlinkAt($ws, [string]$addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            break
        }
    }
}

# Helper: change only the visible display text of an existing hyperlink
# (keeps its target / r:id untouched).
function Set-HyperlinkDisplay($ws, [string]$addr, [string]$text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
            break
        }
    }
}

# Helper: write a plain string value into a cell while reusing the existing
# "True"/"False" shared strings instead of letting them be coerced into
# native booleans, and keeping the resulting style at the cell's previous
# (non hyperlink) style.
function Set-TextValue($ws, [string]$addr, [string]$text) {
    $range = $ws.Range($addr)
    if ($text -eq "") {
        $range.Formula = "'"
    } else {
        $range.Formula = "'" + $text
    }
    $range.Style = "Normal"
}

# New file identities for this handoff run.
$oldFile1 = "7c9669be-7457-4c8c-9b2f-a9069ca0c5d0.md"
$newFile1 = "330fd509-94e5-4011-aec3-2893a3af9f11.md"
$oldFile2 = "addcc77c-c26d-4265-8dcd-e95694c53179.md"
$newFile2 = "ffff5b3b7d53-d086-4552-ab99-5fec2fd5f4ff.md"

$newStatus = "Ready for handoff"
$newHoDate = "2016-08-18 05:02:40"

$newZhXlf = "330fd509-94e5-4011-aec3-2893a3af9f11.edf0d6f0e5e0ca4dce42faf24f1506a81a7a3db3.zh-cn.xlf"
$newDeXlf = "330fd509-94e5-4011-aec3-2893a3af9f11.edf0d6f0e5e0ca4dce42faf24f1506a81a7a3db3.de-de.xlf"
$newZhHandoffDate = "2016-08-18 05:02:35"
$newHandbackDate = "0001-01-01 00:00:00"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $newFile1
$ws1.Range("B2").Value = "e2e\" + $newFile1
Set-HyperlinkDisplay $ws1 '$B$2' ("e2e\" + $newFile1)
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("G2").Value = $newHoDate

$ws1.Range("A3").Value = $newFile2
$ws1.Range("B3").Value = "e2e\" + $newFile2
Set-HyperlinkDisplay $ws1 '$B$3' ("e2e\" + $newFile2)
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws1.Range("G3").Value = $newHoDate

$ws1.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws1.Columns.Item(6).ColumnWidth = 16.333333333333332

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $newFile1
Set-HyperlinkDisplay $ws2 '$A$2' $newFile1
$ws2.Range("G2").Value = $newZhXlf
$ws2.Range("H2").Value = $newZhHandoffDate
Remove-HyperlinkAt $ws2 '$I$2'
Set-TextValue $ws2 "I2" ""
Set-TextValue $ws2 "J2" ""
$ws2.Range("K2").Value = $newHandbackDate

$ws2.Range("A3").Value = $newFile2
Set-HyperlinkDisplay $ws2 '$A$3' $newFile2
Set-TextValue $ws2 "F3" "True"
$ws2.Range("G3").Value = $newZhXlf
Remove-HyperlinkAt $ws2 '$I$3'
Set-TextValue $ws2 "I3" ""
Set-TextValue $ws2 "J3" ""
$ws2.Range("K3").Value = $newHandbackDate

$ws2.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws2.Columns.Item(9).ColumnWidth = 17.833333333333332
$ws2.Columns.Item(10).ColumnWidth = 20.833333333333332

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $newFile1
Set-HyperlinkDisplay $ws3 '$A$2' $newFile1
$ws3.Range("G2").Value = $newDeXlf
$ws3.Range("H2").Value = $newHoDate
Remove-HyperlinkAt $ws3 '$I$2'
Set-TextValue $ws3 "I2" ""
Set-TextValue $ws3 "J2" ""
$ws3.Range("K2").Value = $newHandbackDate

$ws3.Range("A3").Value = $newFile2
Set-HyperlinkDisplay $ws3 '$A$3' $newFile2
Set-TextValue $ws3 "F3" "True"
$ws3.Range("G3").Value = $newDeXlf
Remove-HyperlinkAt $ws3 '$I$3'
Set-TextValue $ws3 "I3" ""
Set-TextValue $ws3 "J3" ""
$ws3.Range("K3").Value = $newHandbackDate

$ws3.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws3.Columns.Item(9).ColumnWidth = 17.833333333333332
$ws3.Columns.Item(10).ColumnWidth = 20.833333333333332
